# Actualización automática 2025-06-02 13:21:56
# Adds a new "PRESUPUESTO" column (G) to the "VENTA MENSUAL" sheet,
# mirroring the formatting of the existing "junio" column (F) and
# filling the new column with 0 for every data/total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy formatting from column F to column G so the new cells pick up the
# same styles (header style, currency number format, total-row style).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("F2:F6").Copy()
$ws.Range("G2:G6").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122)

# New header + values
$ws.Range("G1").Value = "PRESUPUESTO"
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0

# Set column G width (stored width of 17 in the sheet XML).
$ws.Columns.Item(7).ColumnWidth = 16.14
